$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Partner ID"
$ws.Range("J11").Select() | Out-Null
